$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 244.38461
$ws.Range("I6").Value = 120
$ws.Range("K6").Value = 360
$ws.Range("M6").Value = -248
$ws.Range("H29").Value = 4178.3335
$ws.Range("H76").Value = 6094.9
$ws.Range("J76").Value = 6094.9
$ws.Range("L76").Value = 6094.9
$ws.Range("N76").Value = -6724.9
$ws.Range("H79").Value = 6094.9
$ws.Range("J79").Value = 6094.9
$ws.Range("L79").Value = 6094.9
$ws.Range("N79").Value = -8278.9
$ws.Range("H80").Value = 324.86667
$ws.Range("I80").Value = 307.5
$ws.Range("K80").Value = 922.5
$ws.Range("M80").Value = 75.5
$ws.Range("H83").Value = 324.86667
$ws.Range("I83").Value = 307.5
$ws.Range("K83").Value = 2767.5
$ws.Range("M83").Value = 2224.5
$ws.Range("H107").Value = 355.40625
$ws.Range("I107").Value = 361.32144
$ws.Range("K107").Value = 361.32144
$ws.Range("M107").Value = 1558.67856
$ws.Range("H132").Value = 3767.8462
$ws.Range("I132").Value = 3725.7273
$ws.Range("K132").Value = 11177.1819
$ws.Range("M132").Value = -8647.1819
$ws.Range("H138").Value = 9390.937
$ws.Range("J138").Value = 8805.286
$ws.Range("L138").Value = 26415.858
$ws.Range("N138").Value = -36695.858
$ws.Range("H141").Value = 1777
$ws.Range("I141").Value = 1777
$ws.Range("K141").Value = 5331
$ws.Range("M141").Value = -151

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2987.6428
$ws.Range("I2").Value = 2302.4546
$ws.Range("K2").Value = 2302.4546
$ws.Range("M2").Value = -2189.4546
$ws.Range("H24").Value = 16656
$ws.Range("J24").Value = 16656
$ws.Range("L24").Value = 16656
$ws.Range("N24").Value = -17404
$ws.Range("H32").Value = 17570.277
$ws.Range("I32").Value = 11023.042
$ws.Range("K32").Value = 11023.042
$ws.Range("M32").Value = -10736.042
$ws.Range("H74").Value = 2472.5
$ws.Range("I74").Value = 1054.6111
$ws.Range("K74").Value = 1054.6111
$ws.Range("M74").Value = -180.6111000000001
$ws.Range("H77").Value = 2472.5
$ws.Range("I77").Value = 1054.6111
$ws.Range("K77").Value = 5273.0555
$ws.Range("M77").Value = -905.0555000000004
$ws.Range("H100").Value = 16656
$ws.Range("J100").Value = 16656
$ws.Range("L100").Value = 16656
$ws.Range("N100").Value = -18820
$ws.Range("H116").Value = 2987.6428
$ws.Range("I116").Value = 2302.4546
$ws.Range("K116").Value = 2302.4546
$ws.Range("M116").Value = -8.454600000000028
$ws.Range("H132").Value = 2297.55
$ws.Range("I132").Value = 1558.6111
$ws.Range("K132").Value = 4675.8333
$ws.Range("M132").Value = -2145.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2987.6428
$ws.Range("I3").Value = 2302.4546
$ws.Range("K3").Value = 2302.4546
$ws.Range("M3").Value = -2188.4546
$ws.Range("H22").Value = 569.8889
$ws.Range("I22").Value = 597.58826
$ws.Range("K22").Value = 597.58826
$ws.Range("M22").Value = -424.58826
$ws.Range("H94").Value = 2474.8333
$ws.Range("I94").Value = 2376
$ws.Range("K94").Value = 2376
$ws.Range("M94").Value = -1925
$ws.Range("H105").Value = 3717.1614
$ws.Range("I105").Value = 3020.7144
$ws.Range("K105").Value = 3020.7144
$ws.Range("M105").Value = -1273.7144
$ws.Range("H107").Value = 8284.143
$ws.Range("I107").Value = 5795.2
$ws.Range("J107").Value = 14506.5
$ws.Range("K107").Value = 5795.2
$ws.Range("L107").Value = 14506.5
$ws.Range("M107").Value = -3875.2
$ws.Range("N107").Value = -18346.5
$ws.Range("H134").Value = 2874.56
$ws.Range("I134").Value = 1573.5
$ws.Range("J134").Value = 5187.5557
$ws.Range("K134").Value = 4720.5
$ws.Range("L134").Value = 15562.6671
$ws.Range("M134").Value = -2185.5
$ws.Range("N134").Value = -20632.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 125.882355
$ws.Range("I7").Value = 131.3125
$ws.Range("J7").Value = 39
$ws.Range("K7").Value = 131.3125
$ws.Range("L7").Value = 39
$ws.Range("M7").Value = -18.3125
$ws.Range("N7").Value = -265
$ws.Range("H99").Value = 13531.889
$ws.Range("I99").Value = 11664.833
$ws.Range("K99").Value = 11664.833
$ws.Range("M99").Value = -10166.833
$ws.Range("H126").Value = 13531.889
$ws.Range("I126").Value = 11664.833
$ws.Range("K126").Value = 34994.499
$ws.Range("M126").Value = -32524.499
$ws.Range("H132").Value = 2365.3142
$ws.Range("I132").Value = 2215.6072
$ws.Range("K132").Value = 6646.821599999999
$ws.Range("M132").Value = -4116.821599999999
$ws.Range("H134").Value = 2547.2646
$ws.Range("I134").Value = 1939.125
$ws.Range("J134").Value = 4006.8
$ws.Range("K134").Value = 5817.375
$ws.Range("L134").Value = 12020.4
$ws.Range("M134").Value = -3282.375
$ws.Range("N134").Value = -17090.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1655.6364
$ws.Range("I107").Value = 1983
$ws.Range("K107").Value = 5949
$ws.Range("M107").Value = -4029
$ws.Range("H114").Value = 927.44446
$ws.Range("I114").Value = 624.5
$ws.Range("J114").Value = 1533.3334
$ws.Range("K114").Value = 1873.5
$ws.Range("L114").Value = 4600.0002
$ws.Range("M114").Value = 1380.5
$ws.Range("N114").Value = -11108.0002
$ws.Range("H138").Value = 6166.3335
$ws.Range("I138").Value = 1750
$ws.Range("K138").Value = 5250
$ws.Range("M138").Value = -110
$ws.Range("H139").Value = 10385.429
$ws.Range("J139").Value = 9666
$ws.Range("L139").Value = 28998
$ws.Range("N139").Value = -39278
$ws.Range("H140").Value = 3921.2666
$ws.Range("I140").Value = 3447.6155
$ws.Range("K140").Value = 10342.8465
$ws.Range("M140").Value = -5162.8465

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 197.55
$ws.Range("I2").Value = 27.2
$ws.Range("K2").Value = 27.2
$ws.Range("M2").Value = 85.8
$ws.Range("H122").Value = 505893.72
$ws.Range("I122").Value = 106994.3
$ws.Range("K122").Value = 320982.9
$ws.Range("M122").Value = -318532.9
$ws.Range("H132").Value = 4304.7334
$ws.Range("I132").Value = 2889.5
$ws.Range("K132").Value = 8668.5
$ws.Range("M132").Value = -6138.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4966.1665
$ws.Range("I7").Value = 4966.3335
$ws.Range("J7").Value = 4966
$ws.Range("K7").Value = 4966.3335
$ws.Range("L7").Value = 4966
$ws.Range("M7").Value = -4854.3335
$ws.Range("N7").Value = -5190
$ws.Range("H16").Value = 3519.2856
$ws.Range("I16").Value = 3789.1667
$ws.Range("J16").Value = 1900
$ws.Range("K16").Value = 3789.1667
$ws.Range("L16").Value = 1900
$ws.Range("M16").Value = -3619.1667
$ws.Range("N16").Value = -2240
$ws.Range("H61").Value = 8437.888999999999
$ws.Range("I61").Value = 7592.2
$ws.Range("K61").Value = 7592.2
$ws.Range("M61").Value = -7390.2
$ws.Range("H113").Value = 8437.888999999999
$ws.Range("I113").Value = 7592.2
$ws.Range("K113").Value = 7592.2
$ws.Range("M113").Value = -5422.2
$ws.Range("H122").Value = 9099.875
$ws.Range("I122").Value = 7933
$ws.Range("J122").Value = 9800
$ws.Range("K122").Value = 23799
$ws.Range("L122").Value = 29400
$ws.Range("M122").Value = -21349
$ws.Range("N122").Value = -34300
$ws.Range("H126").Value = 4966.1665
$ws.Range("I126").Value = 4966.3335
$ws.Range("J126").Value = 4966
$ws.Range("K126").Value = 14899.0005
$ws.Range("L126").Value = 14898
$ws.Range("M126").Value = -12429.0005
$ws.Range("N126").Value = -19838
$ws.Range("H132").Value = 6500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 19500
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -24560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 10000000
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H81").Value = 3726.2
$ws.Range("I81").Value = 3492.3572
$ws.Range("K81").Value = 6984.7144
$ws.Range("M81").Value = -5923.7144
$ws.Range("H84").Value = 3726.2
$ws.Range("I84").Value = 3492.3572
$ws.Range("K84").Value = 34923.572
$ws.Range("M84").Value = -29619.572
$ws.Range("H132").Value = 2470.4546
$ws.Range("J132").Value = 3285
$ws.Range("L132").Value = 9855
$ws.Range("N132").Value = -14915
$ws.Range("H136").Value = 61029.293
$ws.Range("I136").Value = 1385.3077
$ws.Range("K136").Value = 4155.9231
$ws.Range("M136").Value = -1605.9231
